# major accuracy check update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the polyA isolation protocol label in column G (shared string "NEBNextPoly(A)E7490" -> add "L")
#    All rows share the same string, so update the whole column range.
$ws.Range("G2:G27").Value = "NEBNextPoly(A)E7490L"

# 2. Widen column G to fit the new text
$ws.Columns.Item(7).ColumnWidth = 23.43

# 3. Replace the roboticS1Prep column (I2:I27) FALSE() formulas with literal boolean FALSE values
$ws.Range("I2:I27").Value = $false

# 4. Update the selection to G2:G27 with active cell G2
$ws.Range("G2:G27").Select()
